$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column A, shifting existing data (A-D) to (B-E)
$ws.Range("A1").EntireColumn.Insert()

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# --- Row 2 (values) ---
$ws.Range("A2").Value = "CasesTab"

$ws.Range("B2").Value = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "UNKNOWN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$ws.Range("C2").Value = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "UNKNOWN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

# D2, E2 keep their original content (already shifted by the column insert)

# --- Column widths ---
# Columns B-E keep their inherited widths from the pre-insert A-D columns, so
# only the brand-new column A needs an explicit width (closest value the
# engine's pixel-quantized width model can reach to the target 8.81640625).
$ws.Columns.Item(1).ColumnWidth = 7.98

# --- Row height for row 2 ---
$ws.Rows.Item(2).RowHeight = 174

# --- Styles (wrap text) on query cells B2, C2 ---
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# --- View: selection ---
$ws.Range("C5").Select()
